$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The formulas in B1:D3 ("=1", "=2", "=3", etc.) had an unnecessary
# leading '=' sign. Replace them with plain literal text values
# (force text so the numeric-looking strings are not re-parsed as
# numbers), then restore the default "Normal" style so no formatting
# change lingers on the cells.
function Set-TextValue($range, $text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("B1") "1"
Set-TextValue $ws.Range("C1") "2"
Set-TextValue $ws.Range("D1") "3"

Set-TextValue $ws.Range("B2") "1"
Set-TextValue $ws.Range("C2") "2"
Set-TextValue $ws.Range("D2") "3"

Set-TextValue $ws.Range("B3") "1"
Set-TextValue $ws.Range("C3") "4"
Set-TextValue $ws.Range("D3") "9"
